$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 411, shifting existing rows 411:506 down to 413:508.
$ws.Range("A411:R412").Insert()

# New row 411: Papa / Patagonia / "1a (cosecha lavada)" entry for Macroferia Regional de Talca.
$ws.Range("A411").Value = 5
$ws.Range("B411").Value = "Macroferia Regional de Talca"
$ws.Range("C411").Value = "Maule"
$ws.Range("D411").Value = 44722
$ws.Range("E411").Value = 7
$ws.Range("F411").Value = 100114001
$ws.Range("G411").Value = "Papa"
$ws.Range("H411").Value = "Patagonia"
$ws.Range("I411").Value = "1a (cosecha lavada)"
$ws.Range("J411").Value = 1500
$ws.Range("K411").Value = 6500
$ws.Range("L411").Value = 6500
$ws.Range("M411").Value = 6500
$ws.Range("N411").Value = '$/saco 25 kilos'
$ws.Range("O411").Value = "Región de Los Lagos"
$ws.Range("P411").Value = 260
$ws.Range("Q411").Value = 25
$ws.Range("R411").Value = "Hortaliza"

# New row 412: Papa / Rodeo / "1a (cosecha lavada)" entry for Macroferia Regional de Talca.
$ws.Range("A412").Value = 5
$ws.Range("B412").Value = "Macroferia Regional de Talca"
$ws.Range("C412").Value = "Maule"
$ws.Range("D412").Value = 44722
$ws.Range("E412").Value = 7
$ws.Range("F412").Value = 100114001
$ws.Range("G412").Value = "Papa"
$ws.Range("H412").Value = "Rodeo"
$ws.Range("I412").Value = "1a (cosecha lavada)"
$ws.Range("J412").Value = 1200
$ws.Range("K412").Value = 6500
$ws.Range("L412").Value = 6500
$ws.Range("M412").Value = 6500
$ws.Range("N412").Value = '$/saco 25 kilos'
$ws.Range("O412").Value = "Región de Los Lagos"
$ws.Range("P412").Value = 260
$ws.Range("Q412").Value = 25
$ws.Range("R412").Value = "Hortaliza"
